# The deck originally ships with its custom "Integral" theme living in
# ppt/theme/theme1.xml (the one actually wired to the Slide Master) and the
# generic default "Office Theme" living in ppt/theme/theme2.xml (the one
# wired to the Notes Master). The authored change swaps those two themes'
# contents, so that the Slide Master ends up using the plain "Office" theme
# palette instead of the custom "Integral" one.
#
# Colors are exposed on the PowerPoint object model through
# Master.Theme.ThemeColorScheme, whose 12 slots follow the standard
# MsoThemeColorSchemeIndex ordering:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5-10 accent1..accent6,
#   11 hlink, 12 folHlink
#
# Re-point the Slide Master's theme colors to the stock "Office" palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $scheme.Item($i).RGB = ToRGB $officeColors[$i - 1]
}
